# Update "想去人数" (want-to-go count) figures in the F column across the
# 展览, 演出 and 全部类型 sheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 508
$ws1.Range("F4").Value  = 508
$ws1.Range("F7").Value  = 1384
$ws1.Range("F9").Value  = 460
$ws1.Range("F10").Value = 612
$ws1.Range("F11").Value = 166
$ws1.Range("F16").Value = 1534
$ws1.Range("F26").Value = 711
$ws1.Range("F28").Value = 1375
$ws1.Range("F29").Value = 126

# --- 演出 sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 195

# --- 全部类型 sheet ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 508
$ws4.Range("F5").Value  = 508
$ws4.Range("F8").Value  = 1384
$ws4.Range("F12").Value = 460
$ws4.Range("F13").Value = 612
$ws4.Range("F15").Value = 166
$ws4.Range("F20").Value = 1534
$ws4.Range("F21").Value = 195
$ws4.Range("F38").Value = 711
$ws4.Range("F40").Value = 1375
$ws4.Range("F41").Value = 126
